# whonet_specimen.xlsx edit
# Commit: "added new abx, added new ruling to ESBL per site, Transformed MIC Columns"
#
# The specimen-type lookup table (SPEC_TYPE -> SPEC_ARS code) gains 14 new
# rows (189-202) mapping newly-introduced specimen-type labels onto their
# existing ARS bucket codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (SPEC_TYPE) widens slightly to best-fit the refreshed table -
$ws.Columns.Item(3).AutoFit() | Out-Null

# --- New SPEC_TYPE / SPEC_ARS rows -----------------------------------------
$ws.Cells.Item(189, 3).Value = "et"
$ws.Cells.Item(189, 4).Value = "ta"

$ws.Cells.Item(190, 3).Value = "urine"
$ws.Cells.Item(190, 4).Value = "ur"

$ws.Cells.Item(191, 3).Value = "blood"
$ws.Cells.Item(191, 4).Value = "bl"

$ws.Cells.Item(192, 3).Value = "ure"
$ws.Cells.Item(192, 4).Value = "ur"

$ws.Cells.Item(193, 3).Value = "ot (stone)"
$ws.Cells.Item(193, 4).Value = "ot"

$ws.Cells.Item(194, 3).Value = "ki (stone)"
$ws.Cells.Item(194, 4).Value = "ot"

$ws.Cells.Item(195, 3).Value = "eta"
$ws.Cells.Item(195, 4).Value = "rp"

$ws.Cells.Item(196, 3).Value = "et tip"
$ws.Cells.Item(196, 4).Value = "fb"

$ws.Cells.Item(197, 3).Value = "sputum"
$ws.Cells.Item(197, 4).Value = "rp"

$ws.Cells.Item(198, 3).Value = "cat"
$ws.Cells.Item(198, 4).Value = "fb"

$ws.Cells.Item(199, 3).Value = "b l"
$ws.Cells.Item(199, 4).Value = "bl"

$ws.Cells.Item(200, 3).Value = "csf"
$ws.Cells.Item(200, 4).Value = "sf"

$ws.Cells.Item(201, 3).Value = "sy"
$ws.Cells.Item(201, 4).Value = "fl"

$ws.Cells.Item(202, 3).Value = "pleural fluid"
$ws.Cells.Item(202, 4).Value = "fl"

# --- Move the viewport / selection down to the newly-added rows -----------
$excel.Goto($ws.Range("A175"), $true)
$ws.Range("C192").Select()
